$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row (row 6): add Tool/Comando labels in front of the existing Windows/MAC headers
$ws.Range("D6").Value = "Tool"
$ws.Range("E6").Value = "Comando"

# Re-order / re-populate the keyboard-shortcut table:
#  - Chrome's "Mostar Console" entry now leads (row 7), with no MAC-column value
#  - The Visual Studio Code rows that used to start at row 7 now start at row 8
$ws.Range("D7").Value = "Chrome"
$ws.Range("E7").Value = "Mostar Console"
$ws.Range("F7").Value = "Control + Shift + I"
$ws.Range("G7").ClearContents()

$ws.Range("D8").Value = "Visual Studio Code"
$ws.Range("E8").Value = "Indentar el código"
$ws.Range("F8").Value = "Shift + Alt + F"
$ws.Range("G8").Value = "Shift + Option + F"

$ws.Range("D9").Value = "Visual Studio Code"
$ws.Range("E9").Value = "Aumenta font"
$ws.Range("F9").Value = "Control + [+]"
$ws.Range("G9").Value = "COMMAND + [+]"

$ws.Range("D10").Value = "Visual Studio Code"
$ws.Range("E10").Value = "Disminuye font"
$ws.Range("F10").Value = "Control + [-]"
$ws.Range("G10").Value = "COMMAND + [-]"

# View state: zoomed in on the new Chrome shortcut cell
$excel.ActiveWindow.Zoom = 190
$ws.Range("F7").Select()
